$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

# All rows in the Runmode column (C2:C7) are now "Y" — update every
# cell that previously held "N" (C2, C3, C5, C6, C7); C4 was already "Y".
$ws.Range("C2").Value = "Y"
$ws.Range("C3").Value = "Y"
$ws.Range("C4").Value = "Y"
$ws.Range("C5").Value = "Y"
$ws.Range("C6").Value = "Y"
$ws.Range("C7").Value = "Y"

# Selection moves from the single cell C4 to the whole Runmode range C2:C7.
$ws.Range("C2:C7").Select()
